$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compartments")
$ws.Activate()
$ws.Range("F2:F5").Value = "n"
[void]$ws.Range("F2:F5").Select()
